$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-11-04 01:19:23"

# Insert a new row at position 5 - shifts the existing rows 5 (MT4) and 6
# (Laravel) down to 6 and 7, keeping their content intact.
$ws.Rows.Item(5).Insert()

# New row 5: the newly scraped listing that now tops the list.
$ws.Cells.Item(5,1).Value = $newTimestamp
$ws.Cells.Item(5,2).Value = "社外エンジニア(WEBサイトやシステムのメンテナンス等の保守/改修等)の募集"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5426251"
$ws.Cells.Item(5,7).Value = 53
$ws.Cells.Item(5,8).Value = "◇サイト"

# New row 8 (appended at the end): another newly scraped listing.
$ws.Cells.Item(8,1).Value = $newTimestamp
$ws.Cells.Item(8,2).Value = "〖リモート可〗Delphiエンジニア募集"
$ws.Cells.Item(8,3).Value = "システム開発"
$ws.Cells.Item(8,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(8,5).Value = "期限情報なし"
$ws.Cells.Item(8,6).Value = "https://www.lancers.jp/work/detail/5341051"
$ws.Cells.Item(8,7).Value = 25

# Refresh the "fetched at" timestamp for every row still in the sheet.
$ws.Cells.Item(2,1).Value = $newTimestamp
$ws.Cells.Item(3,1).Value = $newTimestamp
$ws.Cells.Item(4,1).Value = $newTimestamp
$ws.Cells.Item(6,1).Value = $newTimestamp
$ws.Cells.Item(7,1).Value = $newTimestamp

# Rebuild the hyperlinks from scratch: Rows.Insert() does not shift the
# worksheet's hyperlink list, so the old rId->ref mapping would otherwise be
# stale. Clearing and re-adding guarantees F2:F8 all point at the right URL
# again (and each pick up the Hyperlink cell style, matching s="1").
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2,6), "https://www.lancers.jp/work/detail/5405023")
$ws.Hyperlinks.Add($ws.Cells.Item(3,6), "https://www.lancers.jp/work/detail/5251319")
$ws.Hyperlinks.Add($ws.Cells.Item(4,6), "https://www.lancers.jp/work/detail/5425801")
$ws.Hyperlinks.Add($ws.Cells.Item(5,6), "https://www.lancers.jp/work/detail/5426251")
$ws.Hyperlinks.Add($ws.Cells.Item(6,6), "https://www.lancers.jp/work/detail/5426185")
$ws.Hyperlinks.Add($ws.Cells.Item(7,6), "https://www.lancers.jp/work/detail/5426038")
$ws.Hyperlinks.Add($ws.Cells.Item(8,6), "https://www.lancers.jp/work/detail/5341051")
